# Quarterly "dollar_cumulative" income statement update:
#  - drop the oldest reporting period (column D: "6 ماهه منتهی به 1399/06" / "1400-09-01 (5)")
#  - every later period shifts one column to the left
#  - a brand-new period is appended in the now-empty last column (M):
#       "12 ماهه منتهی به 1401/12", published "1402-02-30 (2)"
#  - the old "1401-10-28 (6)" publish-date label is corrected to "1402-02-30 (8)"
#    (this lands on column I after the shift)
#  - new financial figures are filled into the new column M

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete first data column; B..M (incl. D) all shift left by one.
$ws.Range("D1:D28").EntireColumn.Delete()

# Fix the mislabeled publish-date (was "1401-10-28 (6)", now in column I after the shift).
$ws.Range("I9").Value = "1402-02-30 (8)"

# Bring back formatting for the newly-exposed last column (M) by copying column L's look,
# then give it its own (slightly wider) column width like the source file.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$ws.Columns("M").ColumnWidth = 28.17

# Populate the new "12 ماهه منتهی به 1401/12" period column.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-30 (2)"

$ws.Range("M11").Value = 71752
$ws.Range("M12").Value = -37999
$ws.Range("M13").Value = 33753
$ws.Range("M14").Value = -2754
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = -28
$ws.Range("M17").Value = 30972
$ws.Range("M18").Value = -5053
$ws.Range("M19").Value = 752
$ws.Range("M20").Value = 26670
$ws.Range("M21").Value = -4246
$ws.Range("M22").Value = 22424
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 22424
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 6376
$ws.Range("M27").Value = 0
